$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 241; this pushes the existing rows 241..336
# down to 242..337 (matching the diff, where every old row r becomes r+1,
# and a brand-new record lands at row 241).
$ws.Rows.Item(241).Insert()

# Populate the newly inserted row 241 with the new record's data.
# Columns A (Mercado ID) .. I (Calidad) carry forward the same constant
# values used throughout this sheet's block.
$ws.Cells.Item(241, 1).Value = 7
$ws.Cells.Item(241, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(241, 3).Value = 'Ñuble'
$ws.Cells.Item(241, 4).Value = [DateTime]"2022-10-03"
$ws.Cells.Item(241, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(241, 5).Value = 16
$ws.Cells.Item(241, 6).Value = 100114013
$ws.Cells.Item(241, 7).Value = 'Zanahoria'
$ws.Cells.Item(241, 8).Value = 'Sin especificar'
$ws.Cells.Item(241, 9).Value = 'Primera'
$ws.Cells.Item(241, 10).Value = 120
$ws.Cells.Item(241, 11).Value = 10000
$ws.Cells.Item(241, 12).Value = 11000
$ws.Cells.Item(241, 13).Value = 10500
$ws.Cells.Item(241, 14).Value = '$/saco 20 kilos'
$ws.Cells.Item(241, 15).Value = 'Región de Ñuble'
$ws.Cells.Item(241, 16).Value = 525
$ws.Cells.Item(241, 17).Value = 20
$ws.Cells.Item(241, 18).Value = 'Hortaliza'
